{"js": "// Map of old problem text -> new problem text, as described by the diff.\nconst replacements = [\n  [\"74\u00f72=\", \"93\u00f72=\"],\n  [\"78\u00f77=\", \"16\u00f75=\"],\n  [\"54\u00f79=\", \"23\u00f75=\"],\n  [\"45\u00f78=\", \"50\u00f78=\"],\n  [\"43\u00f76=\", \"93\u00f77=\"],\n  [\"61\u00f75=\", \"67\u00f79=\"],\n  [\"59\u00f77=\", \"31\u00f76=\"],\n  [\"60\u00f72=\", \"98\u00f78=\"],\n  [\"33\u00f73=\", \"47\u00f75=\"],\n  [\"95\u00f78=\", \"92\u00f78=\"],\n  [\"46\u00f73=\", \"44\u00f75=\"],\n  [\"79\u00f74=\", \"83\u00f75=\"],\n  [\"72\u00f75=\", \"49\u00f72=\"],\n  [\"30\u00f76=\", \"48\u00f78=\"],\n  [\"41\u00f72=\", \"28\u00f72=\"],\n  [\"73\u00f72=\", \"25\u00f72=\"],\n  [\"99\u00f75=\", \"62\u00f73=\"],\n  [\"20\u00f79=\", \"30\u00f76=\"],\n  [\"82\u00f79=\", \"75\u00f74=\"],\n  [\"30\u00f72=\", \"94\u00f76=\"],\n  [\"14\u00f75=\", \"91\u00f73=\"],\n  [\"65\u00f75=\", \"40\u00f76=\"],\n  [\"80\u00f74=\", \"33\u00f78=\"],\n  [\"31\u00f72=\", \"58\u00f73=\"],\n  [\"75\u00f77=\", \"89\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items,text\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-divided-by-one-digit problem text with its new value,\n# as described by the diff. Every occurrence in the document is unique, so a\n# simple Find/Replace (ReplaceAll) per pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"74\u00f72=\", \"93\u00f72=\"),\n  @(\"78\u00f77=\", \"16\u00f75=\"),\n  @(\"54\u00f79=\", \"23\u00f75=\"),\n  @(\"45\u00f78=\", \"50\u00f78=\"),\n  @(\"43\u00f76=\", \"93\u00f77=\"),\n  @(\"61\u00f75=\", \"67\u00f79=\"),\n  @(\"59\u00f77=\", \"31\u00f76=\"),\n  @(\"60\u00f72=\", \"98\u00f78=\"),\n  @(\"33\u00f73=\", \"47\u00f75=\"),\n  @(\"95\u00f78=\", \"92\u00f78=\"),\n  @(\"46\u00f73=\", \"44\u00f75=\"),\n  @(\"79\u00f74=\", \"83\u00f75=\"),\n  @(\"72\u00f75=\", \"49\u00f72=\"),\n  @(\"30\u00f76=\", \"48\u00f78=\"),\n  @(\"41\u00f72=\", \"28\u00f72=\"),\n  @(\"73\u00f72=\", \"25\u00f72=\"),\n  @(\"99\u00f75=\", \"62\u00f73=\"),\n  @(\"20\u00f79=\", \"30\u00f76=\"),\n  @(\"82\u00f79=\", \"75\u00f74=\"),\n  @(\"30\u00f72=\", \"94\u00f76=\"),\n  @(\"14\u00f75=\", \"91\u00f73=\"),\n  @(\"65\u00f75=\", \"40\u00f76=\"),\n  @(\"80\u00f74=\", \"33\u00f78=\"),\n  @(\"31\u00f72=\", \"58\u00f73=\"),\n  @(\"75\u00f77=\", \"89\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1  # wdFindContinue\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Save()\n"}
